# Workshop3.docx edit: "changed calculator to analyzer"
#
# The only intentional content change in this revision is the word
# "Calculator" -> "Analyzer" inside the highlighted
# "---=== IPC Temperature Calculator ===---" banner text. Word also
# automatically re-homes its "last edit location" bookmark (_GoBack) to
# wherever the edit was made, so we replicate that too.

$d = $word.ActiveDocument

# _GoBack always tracks the most recent edit location in a real Word
# session; remove the stale one (if present) before we make our edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate "Calculator" and replace it with "Analyzer".
$r = $d.Content
$found = $r.Find.Execute("Calculator", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Analyzer", 2)

if ($found) {
    # Drop a transient bookmark right before the replaced run so the
    # engine's run-coalescing pass doesn't fold it into the preceding
    # "---=== IPC Temperature " run (matching the target, where that run
    # boundary is preserved).
    $beforeEdit = $d.Range($r.Start, $r.Start)
    $d.Bookmarks.Add("ZZZ_TEMP_BOUNDARY", $beforeEdit)

    # Re-create _GoBack as a collapsed bookmark immediately after the
    # word we just edited.
    $afterEdit = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $afterEdit)

    # Remove the transient boundary marker now that the surrounding runs
    # have settled.
    $d.Bookmarks("ZZZ_TEMP_BOUNDARY").Delete()
}
